# Update "想去人数" (F5) and "最低票价" (G2, G3) values on both the
# "展览" sheet and the duplicated "全部类型" sheet.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    # G2: 50 -> 58
    $ws.Range("G2").Value = 58

    # G3: 40 -> 50
    $ws.Range("G3").Value = 50

    # F5: 907 -> 909
    $ws.Range("F5").Value = 909
}
